$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Situacao da OS" table: add a tiny table indent (tblInd = 5 dxa)
# -----------------------------------------------------------------
$tSituacao = $d.Tables.Item(4)
$tSituacao.Rows.LeftIndent = 0.25

# -----------------------------------------------------------------
# 2) "Data da Situacao" cell: 06/07/16 -> 13/07/16 (only change the
#    first run's text from "06" to "13", leave the other runs alone)
# -----------------------------------------------------------------
$dataSituacaoCell = $tSituacao.Rows.Item(2).Cells.Item(2)
$rSit = $dataSituacaoCell.Range
$rSit.End = $rSit.Start + 2
$rSit.Text = "13"

# -----------------------------------------------------------------
# 3) "Entrega do Plano" row: fill in "Data Apurada" (13/07/16) and
#    "Documento de Comprovacao" (Acompanhamento da OS)
# -----------------------------------------------------------------
$entregaRow = $tSituacao.Rows.Item(7)

$dataApuradaCell = $entregaRow.Cells.Item(3)
$dataApuradaCell.Range.Text = "13/07/16"

$docComprovacaoCell = $entregaRow.Cells.Item(4)
$docComprovacaoCell.Range.Text = "Acompanhamento da OS"

# -----------------------------------------------------------------
# 4) "Ciclo de Vida da Ordem de Servico" table: remove the stray
#    _GoBack bookmark from the "Chamado Iniciado" paragraph
# -----------------------------------------------------------------
$tCiclo = $d.Tables.Item(5)
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# -----------------------------------------------------------------
# 5) Fill in the new "Chamado Planejado" row (date / event / owner)
# -----------------------------------------------------------------
$novaLinha = $tCiclo.Rows.Item(5)

$novaLinha.Cells.Item(1).Range.Text = "13/07/16"

$eventoCell = $novaLinha.Cells.Item(2)
$eventoCell.Range.ParagraphFormat.Alignment = 1
$eventoCell.Range.Text = "Chamado Planejado"

$novaLinha.Cells.Item(3).Range.Text = "NTConsult " + [char]0x2013 + " Rodrigo Borges"

Write-Output "done"
